$wb = $excel.ActiveWorkbook

# ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 299.5
$ws.Range("I2").Value = 299.5
$ws.Range("K2").Value = 299.5
$ws.Range("M2").Value = -186.5

# ALC row 11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 32.857143
$ws.Range("I11").Value = 32.857143
$ws.Range("K11").Value = 32.857143
$ws.Range("M11").Value = 107.142857

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5249.25

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 973.58826
$ws.Range("I80").Value = 913.8
$ws.Range("K80").Value = 2741.4
$ws.Range("M80").Value = -1743.4

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 973.58826
$ws.Range("I83").Value = 913.8
$ws.Range("K83").Value = 8224.199999999999
$ws.Range("M83").Value = -3232.199999999999

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1808.8
$ws.Range("J100").Value = 2499.75
$ws.Range("L100").Value = 2499.75
$ws.Range("N100").Value = -3581.75

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 6199.75
$ws.Range("I113").Value = 3466.6667
$ws.Range("K113").Value = 3466.6667
$ws.Range("M113").Value = -212.6667000000002

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1177.9231
$ws.Range("I132").Value = 1250.8182
$ws.Range("K132").Value = 3752.4546
$ws.Range("M132").Value = -1222.4546

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2986.2
$ws.Range("I135").Value = 2965.5
$ws.Range("K135").Value = 26689.5
$ws.Range("M135").Value = -24154.5

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3461.75
$ws.Range("I141").Value = 3242
$ws.Range("K141").Value = 9726
$ws.Range("M141").Value = -4546

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2817.2222
$ws.Range("I2").Value = 2856.875
$ws.Range("J2").Value = 2500
$ws.Range("K2").Value = 2856.875
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = -2743.875
$ws.Range("N2").Value = -2726

# ARM row 33
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4590
$ws.Range("I45").Value = 1862.5555
$ws.Range("K45").Value = 1862.5555
$ws.Range("M45").Value = -1485.5555

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7575.3335
$ws.Range("I61").Value = 6579
$ws.Range("J61").Value = 9568
$ws.Range("K61").Value = 6579
$ws.Range("L61").Value = 9568
$ws.Range("M61").Value = -6367
$ws.Range("N61").Value = -9992

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3469.6956
$ws.Range("I102").Value = 2779.158
$ws.Range("K102").Value = 2779.158
$ws.Range("M102").Value = -1157.158

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2817.2222
$ws.Range("I116").Value = 2856.875
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 2856.875
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = -562.875
$ws.Range("N116").Value = -7088

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4754.091
$ws.Range("I132").Value = 3445.6155
$ws.Range("K132").Value = 10336.8465
$ws.Range("M132").Value = -7806.8465

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7575.3335
$ws.Range("I136").Value = 6579
$ws.Range("J136").Value = 9568
$ws.Range("K136").Value = 19737
$ws.Range("L136").Value = 28704
$ws.Range("M136").Value = -17187
$ws.Range("N136").Value = -33804

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2817.2222
$ws.Range("I3").Value = 2856.875
$ws.Range("J3").Value = 2500
$ws.Range("K3").Value = 2856.875
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = -2742.875
$ws.Range("N3").Value = -2728

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 627.4286
$ws.Range("I80").Value = 292.16666
$ws.Range("K80").Value = 292.16666
$ws.Range("M80").Value = 705.83334

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 627.4286
$ws.Range("I83").Value = 292.16666
$ws.Range("K83").Value = 1460.8333
$ws.Range("M83").Value = 3531.1667

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1545.9445
$ws.Range("I86").Value = 1401.9286
$ws.Range("K86").Value = 1401.9286
$ws.Range("M86").Value = -278.9286

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1545.9445
$ws.Range("I89").Value = 1401.9286
$ws.Range("K89").Value = 7009.643
$ws.Range("M89").Value = -1393.643

# BSM row 96
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 17000
$ws.Range("I96").Value = 17000
$ws.Range("K96").Value = 17000
$ws.Range("M96").Value = -14254

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2379.8
$ws.Range("I105").Value = 2379.8
$ws.Range("K105").Value = 2379.8
$ws.Range("M105").Value = -632.8000000000002

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2243.75
$ws.Range("I107").Value = 2000
$ws.Range("J107").Value = 3365
$ws.Range("K107").Value = 2000
$ws.Range("L107").Value = 3365
$ws.Range("M107").Value = -80
$ws.Range("N107").Value = -7205

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5118.4585
$ws.Range("I31").Value = 4076.0667
$ws.Range("K31").Value = 4076.0667
$ws.Range("M31").Value = -3781.0667

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5118.4585
$ws.Range("I34").Value = 4076.0667
$ws.Range("K34").Value = 4076.0667
$ws.Range("M34").Value = -3874.0667

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3149
$ws.Range("I132").Value = 1793.4445
$ws.Range("K132").Value = 5380.333500000001
$ws.Range("M132").Value = -2850.333500000001

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 7491
$ws.Range("I134").Value = 5365.7
$ws.Range("K134").Value = 16097.1
$ws.Range("M134").Value = -13562.1

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 32993778
$ws.Range("I4").Value = 38991984
$ws.Range("K4").Value = 116975952
$ws.Range("M4").Value = -116975840

# CUL row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# CUL row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 9833.333000000001
$ws.Range("J5").Value = 9833.333000000001
$ws.Range("L5").Value = 9833.333000000001
$ws.Range("N5").Value = -10057.333

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3408.2222
$ws.Range("I102").Value = 2903.4285
$ws.Range("K102").Value = 2903.4285
$ws.Range("M102").Value = -1281.4285

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1620
$ws.Range("J122").Value = 2295
$ws.Range("L122").Value = 6885
$ws.Range("N122").Value = -11785

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2622.926
$ws.Range("I126").Value = 1801.375
$ws.Range("J126").Value = 3817.9092
$ws.Range("K126").Value = 5404.125
$ws.Range("L126").Value = 11453.7276
$ws.Range("M126").Value = -2934.125
$ws.Range("N126").Value = -16393.7276

# LTW row 25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 29991
$ws.Range("J25").Value = 29991
$ws.Range("L25").Value = 29991
$ws.Range("N25").Value = -30451

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4492.5
$ws.Range("I122").Value = 4492.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13477.5
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = -11027.5
$ws.Range("M122").ClearContents()

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7394.8184
$ws.Range("I132").Value = 7392.533
$ws.Range("K132").Value = 22177.599
$ws.Range("M132").Value = -19647.599

# WVR row 34
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 1500
$ws.Range("I34").Value = 1500
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1297

# WVR row 70
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 40408.332
$ws.Range("J70").Value = 40408.332
$ws.Range("L70").Value = 40408.332
$ws.Range("N70").Value = -41038.332

# WVR row 73
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 40408.332
$ws.Range("J73").Value = 40408.332
$ws.Range("L73").Value = 40408.332
$ws.Range("N73").Value = -42592.332

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3285.087
$ws.Range("I122").Value = 1860.3636
$ws.Range("J122").Value = 4591.0835
$ws.Range("K122").Value = 5581.0908
$ws.Range("L122").Value = 13773.2505
$ws.Range("M122").Value = -3131.0908
$ws.Range("N122").Value = -18673.2505

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5542.4287
$ws.Range("I136").Value = 3460.4443
$ws.Range("K136").Value = 10381.3329
$ws.Range("M136").Value = -7831.332900000001
